$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new localization rows 54-90 (A = Key, B = Value) ----------
# Cells are written in the exact sequence in which their text values were
# first introduced so the generated shared-strings table lines up with the
# target workbook (uniqueCount 108 -> 179).
$ws.Range("A54").Value = 'stage_1_1_title'
$ws.Range("A55").Value = 'stage_1_1_desc'
$ws.Range("A56").Value = 'stage_1_2_title'
$ws.Range("A57").Value = 'stage_1_2_desc'
$ws.Range("A58").Value = 'stage_1_3_title'
$ws.Range("A59").Value = 'stage_1_3_desc'
$ws.Range("A60").Value = 'stage_2_1_title'
$ws.Range("A61").Value = 'stage_2_1_desc'
$ws.Range("A62").Value = 'stage_2_2_title'
$ws.Range("A63").Value = 'stage_2_2_desc'
$ws.Range("A64").Value = 'stage_2_3_title'
$ws.Range("A65").Value = 'stage_2_3_desc'
$ws.Range("A66").Value = 'stage_3_1_title'
$ws.Range("A67").Value = 'stage_3_1_desc'
$ws.Range("A68").Value = 'stage_3_2_title'
$ws.Range("A69").Value = 'stage_3_2_desc'
$ws.Range("A70").Value = 'stage_3_3_title'
$ws.Range("A71").Value = 'stage_3_3_desc'
$ws.Range("B54").Value = 'Hydrothermal Vent'
$ws.Range("B55").Value = 'A fissure found on the seafloor near volcanically active places. Temperature can reach as high as 867 °F.'
$ws.Range("B56").Value = 'Colon'
$ws.Range("B57").Value = 'The last part of the digestive tract found inside animals. This is where the remaining materials are broken down with the help of bacteria and archaea before excretion.'
$ws.Range("B58").Value = 'Red Sea'
$ws.Range("B59").Value = 'A seawater inlet located between Africa and Asia. The Red Sea is one of the saltiest waters in the world, a perfect place for salt-loving creatures. '
$ws.Range("B60").Value = 'Septic Water'
$ws.Range("B61").Value = 'The kind of water found in sewer where all the wastes accumulate. A perfect home for all sorts of bacteria.'
$ws.Range("B62").Value = 'Pond'
$ws.Range("B65").Value = 'Known as the windpipe that connects the larynx, and the bronchi of the lungs. Many foreign organisms are disposed of here by the hands of white blood cells. Beware of the wandering macrophages.'
$ws.Range("B64").Value = 'Trachea'
$ws.Range("B68").Value = 'Murky Swamp'
$ws.Range("B69").Value = 'A body of freshwater filled with grime. This particular area is filled with toxic bacteria, not ideal for consumption.'
$ws.Range("B70").Value = 'Whirlpool'
$ws.Range("B63").Value = 'A body of freshwater within a land brimming with life. Where there are thriving populations of organisms, so too, will there be predators.'
$ws.Range("B67").Value = 'A body of freshwater within a land brimming with life. A perfect place for predatory organisms.'
$ws.Range("B71").Value = 'Watch out for this downward spiral into the void. However, with this many organisms being pulled in, it’s sure to be a buffet.'
$ws.Range("A78").Value = 'attributeCategoryHazards'
$ws.Range("B78").Value = 'Hazards'
$ws.Range("A79").Value = 'attributeCategoryEnergy'
$ws.Range("B79").Value = 'Energy Sources'
$ws.Range("A81").Value = 'attributeHazardExtremeHighTemperature'
$ws.Range("A82").Value = 'attributeHazardHighSalinity'
$ws.Range("B82").Value = 'High Salinity'
$ws.Range("A83").Value = 'attributeHazardUVRadiation'
$ws.Range("B83").Value = 'UV Radiation'
$ws.Range("A72").Value = 'energySulfur'
$ws.Range("B72").Value = 'Sulfur'
$ws.Range("A73").Value = 'energyH2'
$ws.Range("B73").Value = 'Hydrogen'
$ws.Range("A74").Value = 'energyAminoAcid'
$ws.Range("B74").Value = 'Amino Acid'
$ws.Range("B75").Value = 'Sunlight'
$ws.Range("A75").Value = 'energySunlight'
$ws.Range("A84").Value = 'attributeHazardHighMethane'
$ws.Range("B76").Value = 'Methane'
$ws.Range("A85").Value = 'attributeHazardLowOxygen'
$ws.Range("B85").Value = 'Low Oxygen'
$ws.Range("A76").Value = 'energyMethane'
$ws.Range("A77").Value = 'energyGlucose'
$ws.Range("B77").Value = 'Glucose'
$ws.Range("A80").Value = 'attributeCategoryDanger'
$ws.Range("B80").Value = 'Danger'
$ws.Range("A86").Value = 'attributeDangerHunter'
$ws.Range("A87").Value = 'attributeDangerMacrophage'
$ws.Range("B87").Value = 'Macrophage'
$ws.Range("A88").Value = 'attributeDangerNeutrophil'
$ws.Range("B88").Value = 'Neutrophil'
$ws.Range("B86").Value = 'Stentor'
$ws.Range("A89").Value = 'attributeDangerToxic'
$ws.Range("B89").Value = 'Toxic Bacteria'
$ws.Range("A90").Value = 'attributeDangerWhirlpool'
$ws.Range("B81").Value = 'Extreme Heat'

# Cells whose text duplicates a value entered above (reuse shared string)
$ws.Range("B66").Value = 'Pond'
$ws.Range("B84").Value = 'Methane'
$ws.Range("B90").Value = 'Whirlpool'

# --- Formatting -----------------------------------------------------------
# Longer description cells get vertical-center alignment (new cell style).
$ws.Range("B57").VerticalAlignment = -4108
$ws.Range("B71").VerticalAlignment = -4108
$ws.Range("B72").VerticalAlignment = -4108
$ws.Range("B73").VerticalAlignment = -4108
$ws.Range("B74").VerticalAlignment = -4108
$ws.Range("B75").VerticalAlignment = -4108
$ws.Range("B76").VerticalAlignment = -4108
$ws.Range("B77").VerticalAlignment = -4108
$ws.Range("B80").VerticalAlignment = -4108

# --- View state -------------------------------------------------------------
# Scroll the sheet so the newly added rows are in view and restore the
# active selection that was in place when the workbook was saved.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
$ws.Range("B82").Select()
